$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Adjusted2")
$srcHdr = $src.Range("A1:F1")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "rho"

$headers = @("theta","Jxx","Jyy","beta","gamma","trace_sq")
for ($c=0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c+1).Value = $headers[$c]
}
$hdr = $ws.Range("A1:F1")
$srcHdr.Copy()
$hdr.PasteSpecial(-4122)

$data = @(
    @(0, 0.0007136929823724444, 0.9992864399844273, 0.001386684132653033, 0, 0.9985777442801911),
    @(5, 0.0271634855578083, 0.9728366541305453, -0.162516745045145, 0, 0.9999723954077009),
    @(10, 0.1148364855748853, 0.8851641967665118, -0.3188135790533821, 0, 0.9999872700341491),
    @(15, 0.2489336868765835, 0.7510667813074985, -0.4323925834539165, 0, 0.9999959828974789),
    @(20, 0.4006519051095231, 0.5993485010346, -0.4900289124571287, 0, 0.9999972448481449),
    @(25, 0.5798928673538779, 0.4201071576441872, -0.4935754903954623, 0, 0.9999992909500226),
    @(30, 0.7469526108688559, 0.2530472845091397, -0.4346097085343356, 0, 0.9997423285858503),
    @(35, 0.8749988923054343, 0.1250011303305274, -0.3288897868446408, 0, 0.997585327901073),
    @(40, 0.9635598982079948, 0.03644018136469517, -0.1786957333112875, 0, 0.9936398944598108),
    @(45, 0.9950894009722699, 0.004911038213265518, -0.007874139163917023, 0, 0.9903510383588288),
    @(50, 0.9684668929003352, 0.03153298701405844, 0.160038701472048, 0, 0.9901472238517766),
    @(55, 0.8860815516606284, 0.1139185069612994, 0.3115824374782751, 0, 0.9922851731114052),
    @(60, 0.7598583437729597, 0.240141654066944, 0.4235250518694544, 0, 0.9937996557414414),
    @(65, 0.5938154400749823, 0.4061846162554529, 0.4884871092788996, 0, 0.9948420312173456),
    @(70, 0.4101135175923994, 0.5898867049572349, 0.4891098270750975, 0, 0.9946162678801789),
    @(75, 0.2627534397818948, 0.7372465756407056, 0.4373131198750824, 0, 0.9950574130409209),
    @(80, 0.1273358720628602, 0.8726641438642255, 0.3306994357323851, 0, 0.9964813658877267),
    @(85, 0.03864335819542286, 0.9613566552656744, 0.1889503866152527, 0, 0.9971044249603314),
    @(90, 0.001059429823470423, 0.998940570739852, 0.01119670036113816, 0, 0.9981341184595662)
)

for ($r=0; $r -lt $data.Length; $r++) {
    for ($c=0; $c -lt $data[$r].Length; $c++) {
        $ws.Cells.Item($r+2, $c+1).Value = $data[$r][$c]
    }
}
Write-Output "done"
